$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 1215
$ws.Cells.Item(4, 6).Value = 1276
$ws.Cells.Item(7, 6).Value = 552
$ws.Cells.Item(8, 6).Value = 16
$ws.Cells.Item(9, 6).Value = 343
$ws.Cells.Item(10, 6).Value = 55
$ws.Cells.Item(11, 6).Value = 1268
$ws.Cells.Item(12, 6).Value = 29135
$ws.Cells.Item(13, 6).Value = 4231
$ws.Cells.Item(15, 6).Value = 263
$ws.Cells.Item(16, 6).Value = 489
$ws.Cells.Item(17, 6).Value = 43
$ws.Cells.Item(19, 6).Value = 14
$ws.Cells.Item(21, 6).Value = 344
$ws.Cells.Item(22, 6).Value = 634
$ws.Cells.Item(23, 6).Value = 277
$ws.Cells.Item(24, 6).Value = 282
$ws.Cells.Item(25, 6).Value = 358
$ws.Cells.Item(29, 6).Value = 666
$ws.Cells.Item(32, 6).Value = 549
$ws.Cells.Item(35, 6).Value = 643
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(4, 6).Value = 19
$ws.Cells.Item(6, 6).Value = 386
$ws.Cells.Item(7, 6).Value = 892
$ws.Cells.Item(11, 6).Value = 275
$ws.Cells.Item(12, 6).Value = 4251
$ws.Cells.Item(14, 6).Value = 190
$ws.Cells.Item(21, 6).Value = 149
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 299
$ws.Cells.Item(3, 6).Value = 263
$ws.Cells.Item(4, 6).Value = 1205
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 299
$ws.Cells.Item(3, 6).Value = 263
$ws.Cells.Item(4, 6).Value = 1205
$ws.Cells.Item(5, 6).Value = 5
$ws.Cells.Item(6, 6).Value = 19
$ws.Cells.Item(7, 6).Value = 386
$ws.Cells.Item(9, 6).Value = 892
$ws.Cells.Item(10, 6).Value = 1215
$ws.Cells.Item(11, 6).Value = 1276
$ws.Cells.Item(13, 6).Value = 552
$ws.Cells.Item(14, 6).Value = 16
$ws.Cells.Item(15, 6).Value = 343
$ws.Cells.Item(17, 6).Value = 55
$ws.Cells.Item(18, 6).Value = 1268
$ws.Cells.Item(21, 6).Value = 275
$ws.Cells.Item(23, 6).Value = 190
$ws.Cells.Item(28, 6).Value = 489
$ws.Cells.Item(29, 6).Value = 43
$ws.Cells.Item(30, 6).Value = 14
$ws.Cells.Item(34, 6).Value = 344
$ws.Cells.Item(35, 6).Value = 634
$ws.Cells.Item(36, 6).Value = 277
$ws.Cells.Item(40, 6).Value = 666
$ws.Cells.Item(44, 6).Value = 149
$ws.Cells.Item(48, 6).Value = 643
